$d = $word.ActiveDocument

# Budget resume (BordereauCnas) monthly refresh: Aout -> Septembre, plus the
# corresponding headcount/amount figures and the spelled-out total.

# 1. "MOIS DU Aout 2020" -> "MOIS DU Septembre 2020"
$r1 = $d.Content.Find.Execute("Aout", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Septembre", 2)

# 2. GEST count 4194 -> 4156
$r2 = $d.Content.Find.Execute("4194", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4156", 2)

# 3. MONT SOUMIS AUX COTISA total 75 492 000,00 -> 74 808 000,00
$r3 = $d.Content.Find.Execute("75 492 000,00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "74 808 000,00", 2)

# 4. Both "3 774 600,00" occurrences (5% line + recap line) -> "3 740 400,00"
#    (wdReplaceAll / Replace:=2 replaces every match in one call)
$r4 = $d.Content.Find.Execute("3 774 600,00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3 740 400,00", 2)

# 5. Amount-in-words sentence update
$r5 = $d.Content.Find.Execute("TROIS MILLIONS SEPT CENT SOIXANTE-QUATORZE MILLE SIX CENTS  ",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "TROIS MILLIONS SEPT CENT QUARANTE MILLE QUATRE CENTS  ", 2)

Write-Output "Aout->Septembre: $r1; 4194->4156: $r2; 75492->74808: $r3; 3774600->3740400: $r4; words: $r5"
